# Weekly crime-data refresh for 68th Precinct CompStat report.
# Advances the reporting week (volume/date header) and updates all
# Week-to-Date / 28-Day / Year-to-Date / 2-Year complaint figures and
# percent-change columns for rows 15-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (rich-text shared strings): bump volume/issue number and
# --- roll the reporting week forward by one week. Characters() edits the
# --- specific run in place so only the targeted substring changes.
$ws.Range("A8").Characters(21, 2).Text = "34"
$ws.Range("C9").Characters(27, 9).Text = "8/19/2024"
$ws.Range("C9").Characters(47, 9).Text = "8/25/2024"

# --- Cells that flip between a numeric value and the "no data" text
# --- placeholders ("0" / "***.*"). Copy from a stable same-shaped cell so
# --- the destination picks up both the correct value/type and formatting.
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("E14").Copy($ws.Range("E20"))

$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("C33").Value = 1
$ws.Range("F33").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0

# --- Remaining cells: same type before/after, only the number changes.

$ws.Range("N15").Value = -77.777777777777
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -66.666666666666
$ws.Range("J16").Value = 48
$ws.Range("K16").Value = -12.5
$ws.Range("L16").Value = 27.272727272727
$ws.Range("M16").Value = -38.235294117647
$ws.Range("N16").Value = -88.429752066115
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = 36.363636363636
$ws.Range("I17").Value = 82
$ws.Range("J17").Value = 95
$ws.Range("K17").Value = -13.684210526315
$ws.Range("L17").Value = 12.328767123287
$ws.Range("M17").Value = 24.242424242424
$ws.Range("N17").Value = -55.913978494623
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 54
$ws.Range("K18").Value = -18.181818181818
$ws.Range("L18").Value = -1.818181818181
$ws.Range("M18").Value = -66.037735849056
$ws.Range("N18").Value = -92.329545454545
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -52.777777777777
$ws.Range("I19").Value = 243
$ws.Range("J19").Value = 286
$ws.Range("K19").Value = -15.034965034965
$ws.Range("L19").Value = -30.769230769230
$ws.Range("M19").Value = 19.704433497536
$ws.Range("N19").Value = -18.456375838926
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 15.384615384615
$ws.Range("I20").Value = 135
$ws.Range("K20").Value = 68.75
$ws.Range("L20").Value = 80
$ws.Range("M20").Value = 26.168224299065
$ws.Range("N20").Value = -89.436619718309
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 23.076923076923
$ws.Range("F21").Value = 54
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = -23.943661971831
$ws.Range("I21").Value = 561
$ws.Range("J21").Value = 583
$ws.Range("K21").Value = -3.773584905660
$ws.Range("L21").Value = -5.872483221476
$ws.Range("M21").Value = -8.032786885245
$ws.Range("N21").Value = -80.315789473684
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -13.513513513513
$ws.Range("F24").Value = 145
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = 28.318584070796
$ws.Range("I24").Value = 946
$ws.Range("J24").Value = 995
$ws.Range("K24").Value = -4.924623115577
$ws.Range("L24").Value = -21.100917431192
$ws.Range("M24").Value = 27.493261455525
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -23.809523809523
$ws.Range("F25").Value = 84
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = 44.827586206896
$ws.Range("I25").Value = 598
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 19.6
$ws.Range("L25").Value = -15.297450424929
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 23.076923076923
$ws.Range("I26").Value = 235
$ws.Range("J26").Value = 243
$ws.Range("K26").Value = -3.292181069958
$ws.Range("L26").Value = 12.440191387559
$ws.Range("M26").Value = -1.260504201680
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 36
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 80
$ws.Range("L28").Value = 9.090909090909
$ws.Range("I33").Value = 6
$ws.Range("K33").Value = 20
$ws.Range("L33").Value = 500
